$wb = $excel.ActiveWorkbook

# Update the "CreateProject" sheet: C2 value changes from "aaa" to "abc"
$createProject = $wb.Worksheets.Item("CreateProject")
$createProject.Range("C2").Value = "abc"

# Make CreateProject the active sheet and update its selection to C3
$createProject.Activate()
$createProject.Range("C3").Select()
